# Update the lattice-multiplication exercise table: each of the 15 cells
# (5 rows x 3 columns) keeps its layout (header / multiplicand split /
# "----" / two lattice-left digits) but gets new multiplication problems.
#
# New "A x B" values, in row-major order (row 1 col 1, row 1 col 2, ...):
$newProblems = @(
    "47 x 82", "67 x 77", "28 x 59",
    "37 x 89", "66 x 86", "75 x 56",
    "77 x 94", "42 x 58", "66 x 67",
    "27 x 96", "74 x 46", "58 x 20",
    "26 x 44", "90 x 42", "64 x 44"
)

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$nl = [char]11   # manual line break -> <w:br/> when written back to OOXML
$i = 0

for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    $row = $tbl.Rows.Item($r)
    for ($c = 1; $c -le $row.Cells.Count; $c++) {
        $cell = $row.Cells.Item($c)

        $problem = $newProblems[$i]
        $i = $i + 1

        $parts = $problem.Split(" x ")
        # Split(" x ") on a single-char delimiter set also splits on the
        # spaces themselves, so filter out the empty pieces it leaves behind.
        $nums = @()
        foreach ($p in $parts) { if ($p -ne "") { $nums += $p } }
        $a = $nums[0]
        $b = $nums[1]

        # Use the -f format operator throughout: plain "+" concatenation of
        # strings that happen to look numeric (e.g. "  " + "8" + "    " + "2")
        # gets silently coerced to numeric addition by this interpreter, so
        # string interpolation via -f sidesteps that entirely.
        $line1 = $problem
        $line2 = "  {0}    {1}" -f $b.Substring(0,1), $b.Substring(1,1)
        $line3 = "  ----"
        $line4 = "{0}|    |" -f $a.Substring(0,1)
        $line5 = "{0}|    |" -f $a.Substring(1,1)

        $newText = "{0}{5}{1}{5}{2}{5}{3}{5}{4}" -f $line1, $line2, $line3, $line4, $line5, $nl

        $cell.Range.Text = $newText
    }
}

Write-Output "Updated $i cells"
